# fix(publipostage): Refactor synthetic array
#
# Inserts a new "statut_name" column right after "statut_label" (i.e. before
# the former column C / "NCTId"), pushing every subsequent column one slot to
# the right (NCTId: C->D, eudraCT: D->E, CTIS: E->F, completion_year: F->G,
# clinical_trial_title: G->H, acronym: H->I, results_1y: I->J,
# results_3y: J->K, results: K->L, intervention_type: L->M).
#
# The new column is populated with a human readable label derived from the
# existing "statut_label" column (B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49

# Insert a brand-new column before the current column C ("NCTId"). This
# shifts C..L to D..M automatically (values, formats, everything) while
# leaving a blank column C ready to receive the new "statut_name" data.
$ws.Columns.Item(3).Insert()

# Header for the freshly inserted column (copy the neighbouring header's
# formatting - bold font, centered alignment, border - onto it).
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 3).Value = "statut_name"

# Map each existing "statut_label" (column B) value to its long-form
# "statut_name" text.
$labelToName = @{
    "noir"   = "pas de résultat ni de publication";
    "orange" = "résultat et / ou publication posté dans les 36 mois";
    "vert"   = "résultat et / ou publication posté dans les 12 mois";
    "rouge"  = "résultat et / ou publication posté";
}

for ($row = 2; $row -le $lastRow; $row++) {
    $label = $ws.Cells.Item($row, 2).Value2
    $name = $labelToName[$label]
    $ws.Cells.Item($row, 3).Value = $name
}
